$d = $word.ActiveDocument

# 1) Insert a new paragraph with the repo link right after the first
#    (header) paragraph, before the paragraph holding the first picture.
$introPara = $d.Paragraphs.Item(1)
$r = $introPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$linkPara = $d.Paragraphs.Item(2)
$linkPara.Range.Text = "https://github.com/lap202/csd-340"

# 2) Mark the runs that hold the three screenshots as NoProof (adds
#    <w:noProof/> to their <w:rPr>), matching Word's normal behavior for
#    inline pictures.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.NoProofing = 1
    }
}
